$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E12")

# Remove the placeholder picture from this sheet
if ($ws.Shapes.Count -gt 0) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        [void]$ws.Shapes.Item($i).Delete()
    }
}

# Header row (row 1): rename Time header, add Upper Limit header
$ws.Range("A1").Value = "time (s)"
$ws.Range("B1").Value = "Thrust (N)"
$ws.Range("D1").Value = "Slope (N/s)"
$ws.Range("E1").Value = "Upper Limit"

# Data rows (Time/Thrust plus derived Slope/Upper Time Bound columns)
$ws.Cells.Item(2,1).Value = 0.002
$ws.Cells.Item(2,2).Value = 1.67799999999999994
$ws.Cells.Item(3,1).Value = 0.05
$ws.Cells.Item(3,2).Value = 4.45899999999999963
$ws.Cells.Item(3,4).Value = 57.93749999999999289
$ws.Cells.Item(3,5).Value = 0.05
$ws.Cells.Item(4,1).Value = 0.10000000000000001
$ws.Cells.Item(4,2).Value = 10.43099999999999916
$ws.Cells.Item(4,4).Value = 119.43999999999998352
$ws.Cells.Item(4,5).Value = 0.10000000000000001
$ws.Cells.Item(5,1).Value = 0.20000000000000001
$ws.Cells.Item(5,2).Value = 24.15200000000000102
$ws.Cells.Item(5,4).Value = 137.21000000000000796
$ws.Cells.Item(5,5).Value = 0.20000000000000001
$ws.Cells.Item(6,1).Value = 0.27400000000000002
$ws.Cells.Item(6,2).Value = 31.95899999999999963
$ws.Cells.Item(6,4).Value = 105.49999999999997158
$ws.Cells.Item(6,5).Value = 0.27400000000000002
$ws.Cells.Item(7,1).Value = 0.29199999999999998
$ws.Cells.Item(7,2).Value = 32.70199999999999818
$ws.Cells.Item(7,4).Value = 41.27777777777778567
$ws.Cells.Item(7,5).Value = 0.29199999999999998
$ws.Cells.Item(8,1).Value = 0.31
$ws.Cells.Item(8,2).Value = 27.95700000000000074
$ws.Cells.Item(8,4).Value = -263.61111111111074479
$ws.Cells.Item(8,5).Value = 0.31
$ws.Cells.Item(9,1).Value = 0.32000000000000001
$ws.Cells.Item(9,2).Value = 25.95700000000000074
$ws.Cells.Item(9,4).Value = -199.99999999999982947
$ws.Cells.Item(9,5).Value = 0.32000000000000001
$ws.Cells.Item(10,1).Value = 0.33000000000000002
$ws.Cells.Item(10,2).Value = 22.85699999999999932
$ws.Cells.Item(10,4).Value = -309.99999999999988631
$ws.Cells.Item(10,5).Value = 0.33000000000000002
$ws.Cells.Item(11,1).Value = 0.34000000000000002
$ws.Cells.Item(11,2).Value = 19.12800000000000011
$ws.Cells.Item(11,4).Value = -372.89999999999957936
$ws.Cells.Item(11,5).Value = 0.34000000000000002
$ws.Cells.Item(12,1).Value = 0.34999999999999998
$ws.Cells.Item(12,2).Value = 16.45499999999999829
$ws.Cells.Item(12,4).Value = -267.30000000000143245
$ws.Cells.Item(12,5).Value = 0.34999999999999998
$ws.Cells.Item(13,1).Value = 0.35999999999999999
$ws.Cells.Item(13,2).Value = 15.31400000000000006
$ws.Cells.Item(13,4).Value = -114.09999999999972431
$ws.Cells.Item(13,5).Value = 0.35999999999999999
$ws.Cells.Item(14,1).Value = 0.38
$ws.Cells.Item(14,2).Value = 13.85299999999999976
$ws.Cells.Item(14,4).Value = -73.04999999999995453
$ws.Cells.Item(14,5).Value = 0.38
$ws.Cells.Item(15,1).Value = 0.39000000000000001
$ws.Cells.Item(15,2).Value = 13.43599999999999994
$ws.Cells.Item(15,4).Value = -41.699999999999946
$ws.Cells.Item(15,5).Value = 0.39000000000000001
$ws.Cells.Item(16,1).Value = 0.40000000000000002
$ws.Cells.Item(16,2).Value = 13.2710000000000008
$ws.Cells.Item(16,4).Value = -16.49999999999990052
$ws.Cells.Item(16,5).Value = 0.40000000000000002
$ws.Cells.Item(17,1).Value = 0.45000000000000001
$ws.Cells.Item(17,2).Value = 12.07000000000000028
$ws.Cells.Item(17,4).Value = -24.02000000000001734
$ws.Cells.Item(17,5).Value = 0.45000000000000001
$ws.Cells.Item(18,1).Value = 0.5
$ws.Cells.Item(18,2).Value = 11.52200000000000024
$ws.Cells.Item(18,4).Value = -10.96000000000000263
$ws.Cells.Item(18,5).Value = 0.5
$ws.Cells.Item(19,1).Value = 0.55000000000000004
$ws.Cells.Item(19,2).Value = 11.26600000000000001
$ws.Cells.Item(19,4).Value = -5.12000000000000011
$ws.Cells.Item(19,5).Value = 0.55000000000000004
$ws.Cells.Item(20,1).Value = 0.59999999999999998
$ws.Cells.Item(20,2).Value = 10.73600000000000065
$ws.Cells.Item(20,4).Value = -10.60000000000000142
$ws.Cells.Item(20,5).Value = 0.59999999999999998
$ws.Cells.Item(21,1).Value = 0.65000000000000002
$ws.Cells.Item(21,2).Value = 10.77699999999999925
$ws.Cells.Item(21,4).Value = 0.81999999999997109
$ws.Cells.Item(21,5).Value = 0.65000000000000002
$ws.Cells.Item(22,1).Value = 0.69999999999999996
$ws.Cells.Item(22,2).Value = 10.2759999999999998
$ws.Cells.Item(22,4).Value = -10.02000000000000313
$ws.Cells.Item(22,5).Value = 0.69999999999999996
$ws.Cells.Item(23,1).Value = 0.80000000000000004
$ws.Cells.Item(23,2).Value = 10.10500000000000043
$ws.Cells.Item(23,4).Value = -1.70999999999999219
$ws.Cells.Item(23,5).Value = 0.80000000000000004
$ws.Cells.Item(24,1).Value = 0.90000000000000002
$ws.Cells.Item(24,2).Value = 9.91999999999999993
$ws.Cells.Item(24,4).Value = -1.85000000000000542
$ws.Cells.Item(24,5).Value = 0.90000000000000002
$ws.Cells.Item(25,1).Value = 1.0
$ws.Cells.Item(25,2).Value = 9.69299999999999962
$ws.Cells.Item(25,4).Value = -2.27000000000000357
$ws.Cells.Item(25,5).Value = 1.0
$ws.Cells.Item(26,1).Value = 1.31000000000000005
$ws.Cells.Item(26,2).Value = 9.75900000000000034
$ws.Cells.Item(26,4).Value = 0.21290322580645391
$ws.Cells.Item(26,5).Value = 1.31000000000000005
$ws.Cells.Item(27,1).Value = 1.31600000000000006
$ws.Cells.Item(27,2).Value = 10.69599999999999973
$ws.Cells.Item(27,4).Value = 156.16666666666642982
$ws.Cells.Item(27,5).Value = 1.31600000000000006
$ws.Cells.Item(28,1).Value = 1.33000000000000007
$ws.Cells.Item(28,2).Value = 9.62800000000000011
$ws.Cells.Item(28,4).Value = -76.28571428571419233
$ws.Cells.Item(28,5).Value = 1.33000000000000007
$ws.Cells.Item(29,1).Value = 2.37999999999999989
$ws.Cells.Item(29,2).Value = 9.86999999999999922
$ws.Cells.Item(29,4).Value = 0.23047619047618967
$ws.Cells.Item(29,5).Value = 2.37999999999999989
$ws.Cells.Item(30,1).Value = 2.39999999999999991
$ws.Cells.Item(30,2).Value = 6.44200000000000017
$ws.Cells.Item(30,4).Value = -171.39999999999980673
$ws.Cells.Item(30,5).Value = 2.39999999999999991
$ws.Cells.Item(31,1).Value = 2.41999999999999993
$ws.Cells.Item(31,2).Value = 3.67399999999999993
$ws.Cells.Item(31,4).Value = -138.399999999999892
$ws.Cells.Item(31,5).Value = 2.41999999999999993
$ws.Cells.Item(32,1).Value = 2.43999999999999995
$ws.Cells.Item(32,2).Value = 0.0
$ws.Cells.Item(32,4).Value = -183.69999999999984652
$ws.Cells.Item(32,5).Value = 2.43999999999999995

# Column widths
$ws.Columns.Item(2).ColumnWidth = 9.6
$ws.Columns.Item(4).ColumnWidth = 12.14

# Selection matching the saved view state
[void]$ws.Range("D3:E32").Select()
